$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 277
$ws1.Range("F4").Value = 931

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 277
$ws4.Range("F5").Value = 931
